# Apply targeted value updates as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E3").Value = 12.24399999999999
$ws.Range("B9").Value = 8.594800000000005
$ws.Range("B18").Value = 4.478700000000005
$ws.Range("B20").Value = 5.704
